$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "result" column header and computed values
$ws.Range("D1").Value = "result"
$ws.Range("D2").Value = 8
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 27
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 4

# Update the active selection to match the post-edit state
$ws.Range("F18").Select()
